$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 1787.8611  # H17: 1788.0278 -> 1787.8611
$ws.Cells.Item(17, 10).Value = 2160.1428  # J17: 2160.4285 -> 2160.1428
$ws.Cells.Item(17, 12).Value = 6480.428400000001  # L17: 6481.2855 -> 6480.428400000001
$ws.Cells.Item(17, 14).Value = -6816.428400000001  # N17: -6817.2855 -> -6816.428400000001

$ws.Cells.Item(87, 8).Value = 53514.5  # H87: 54150.445 -> 53514.5
$ws.Cells.Item(87, 10).Value = 59445.145  # J87: 59419.25 -> 59445.145
$ws.Cells.Item(87, 12).Value = 59445.145  # L87: 59419.25 -> 59445.145
$ws.Cells.Item(87, 14).Value = -61941.145  # N87: -61915.25 -> -61941.145

$ws.Cells.Item(90, 8).Value = 53514.5  # H90: 54150.445 -> 53514.5
$ws.Cells.Item(90, 10).Value = 59445.145  # J90: 59419.25 -> 59445.145
$ws.Cells.Item(90, 12).Value = 178335.435  # L90: 178257.75 -> 178335.435
$ws.Cells.Item(90, 14).Value = -190815.435  # N90: -190737.75 -> -190815.435

$ws.Cells.Item(92, 8).Value = 828  # H92: 889.1667 -> 828
$ws.Cells.Item(92, 9).Value = 158  # I92: 167.14285 -> 158
$ws.Cells.Item(92, 11).Value = 158  # K92: 167.14285 -> 158
$ws.Cells.Item(92, 13).Value = 1090  # M92: 1080.85715 -> 1090

$ws.Cells.Item(98, 8).Value = 1256.8948  # H98: 1290.4445 -> 1256.8948
$ws.Cells.Item(98, 9).Value = 1256.8948  # I98: 1290.4445 -> 1256.8948
$ws.Cells.Item(98, 11).Value = 1256.8948  # K98: 1290.4445 -> 1256.8948
$ws.Cells.Item(98, 13).Value = 241.1052  # M98: 207.5554999999999 -> 241.1052

$ws.Cells.Item(112, 8).Value = 2833.3333  # H112: 2750 -> 2833.3333
$ws.Cells.Item(112, 10).Value = 3250  # J112: 3000 -> 3250
$ws.Cells.Item(112, 12).Value = 9750  # L112: 9000 -> 9750
$ws.Cells.Item(112, 14).Value = -11966  # N112: -11216 -> -11966

$ws.Cells.Item(122, 8).Value = 1256.8948  # H122: 1290.4445 -> 1256.8948
$ws.Cells.Item(122, 9).Value = 1256.8948  # I122: 1290.4445 -> 1256.8948
$ws.Cells.Item(122, 11).Value = 3770.6844  # K122: 3871.3335 -> 3770.6844
$ws.Cells.Item(122, 13).Value = -1320.6844  # M122: -1421.3335 -> -1320.6844

$ws.Cells.Item(125, 8).Value = 4112.25  # H125: 4214 -> 4112.25
$ws.Cells.Item(125, 9).Value = 3816.3333  # I125: 3899.6 -> 3816.3333
$ws.Cells.Item(125, 11).Value = 34346.9997  # K125: 35096.4 -> 34346.9997
$ws.Cells.Item(125, 13).Value = -31886.9997  # M125: -32636.4 -> -31886.9997

$ws.Cells.Item(131, 8).Value = 1296.6666  # H131: 1411.25 -> 1296.6666
$ws.Cells.Item(131, 9).Value = 674.6  # I131: 748.25 -> 674.6
$ws.Cells.Item(131, 11).Value = 2023.8  # K131: 2244.75 -> 2023.8
$ws.Cells.Item(131, 13).Value = 3016.2  # M131: 2795.25 -> 3016.2

$ws.Cells.Item(132, 8).Value = 14244.263  # H132: 14257.421 -> 14244.263
$ws.Cells.Item(132, 9).Value = 16380.077  # I132: 16399.309 -> 16380.077
$ws.Cells.Item(132, 11).Value = 49140.231  # K132: 49197.927 -> 49140.231
$ws.Cells.Item(132, 13).Value = -46610.231  # M132: -46667.927 -> -46610.231

$ws.Cells.Item(137, 8).Value = 2671.16  # H137: 2159.4707 -> 2671.16
$ws.Cells.Item(137, 9).Value = 1959.2858  # I137: 1615.1578 -> 1959.2858
$ws.Cells.Item(137, 10).Value = 3577.182  # J137: 2848.9333 -> 3577.182
$ws.Cells.Item(137, 11).Value = 5877.857400000001  # K137: 4845.4734 -> 5877.857400000001
$ws.Cells.Item(137, 12).Value = 10731.546  # L137: 8546.7999 -> 10731.546
$ws.Cells.Item(137, 13).Value = -3327.857400000001  # M137: -2295.4734 -> -3327.857400000001
$ws.Cells.Item(137, 14).Value = -15831.546  # N137: -13646.7999 -> -15831.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2398.2307  # H2: 2225.3635 -> 2398.2307
$ws.Cells.Item(2, 9).Value = 2017  # I2: 1721 -> 2017
$ws.Cells.Item(2, 11).Value = 2017  # K2: 1721 -> 2017
$ws.Cells.Item(2, 13).Value = -1904  # M2: -1608 -> -1904

$ws.Cells.Item(4, 8).Value = 156.5  # H4: 181.5 -> 156.5
$ws.Cells.Item(4, 9).Value = 156.5  # I4: 181.5 -> 156.5
$ws.Cells.Item(4, 11).Value = 156.5  # K4: 181.5 -> 156.5
$ws.Cells.Item(4, 13).Value = -40.5  # M4: -65.5 -> -40.5

$ws.Cells.Item(74, 8).Value = 4111.9546  # H74: 4123.591 -> 4111.9546
$ws.Cells.Item(74, 9).Value = 3103.2666  # I74: 3293.2144 -> 3103.2666
$ws.Cells.Item(74, 10).Value = 6273.4287  # J74: 5576.75 -> 6273.4287
$ws.Cells.Item(74, 11).Value = 3103.2666  # K74: 3293.2144 -> 3103.2666
$ws.Cells.Item(74, 12).Value = 6273.4287  # L74: 5576.75 -> 6273.4287
$ws.Cells.Item(74, 13).Value = -2229.2666  # M74: -2419.2144 -> -2229.2666
$ws.Cells.Item(74, 14).Value = -8021.4287  # N74: -7324.75 -> -8021.4287

$ws.Cells.Item(77, 8).Value = 4111.9546  # H77: 4123.591 -> 4111.9546
$ws.Cells.Item(77, 9).Value = 3103.2666  # I77: 3293.2144 -> 3103.2666
$ws.Cells.Item(77, 10).Value = 6273.4287  # J77: 5576.75 -> 6273.4287
$ws.Cells.Item(77, 11).Value = 15516.333  # K77: 16466.072 -> 15516.333
$ws.Cells.Item(77, 12).Value = 31367.1435  # L77: 27883.75 -> 31367.1435
$ws.Cells.Item(77, 13).Value = -11148.333  # M77: -12098.072 -> -11148.333
$ws.Cells.Item(77, 14).Value = -40103.14350000001  # N77: -36619.75 -> -40103.14350000001

$ws.Cells.Item(102, 8).Value = 5932.65  # H102: 6142.1665 -> 5932.65
$ws.Cells.Item(102, 9).Value = 4710.2666  # I102: 4812.3076 -> 4710.2666
$ws.Cells.Item(102, 11).Value = 4710.2666  # K102: 4812.3076 -> 4710.2666
$ws.Cells.Item(102, 13).Value = -3088.2666  # M102: -3190.3076 -> -3088.2666

$ws.Cells.Item(116, 8).Value = 2398.2307  # H116: 2225.3635 -> 2398.2307
$ws.Cells.Item(116, 9).Value = 2017  # I116: 1721 -> 2017
$ws.Cells.Item(116, 11).Value = 2017  # K116: 1721 -> 2017
$ws.Cells.Item(116, 13).Value = 277  # M116: 573 -> 277

$ws.Cells.Item(132, 8).Value = 2741.923  # H132: 2887.1667 -> 2741.923
$ws.Cells.Item(132, 9).Value = 1421.7778  # I132: 1474.625 -> 1421.7778
$ws.Cells.Item(132, 11).Value = 4265.3334  # K132: 4423.875 -> 4265.3334
$ws.Cells.Item(132, 13).Value = -1735.3334  # M132: -1893.875 -> -1735.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2398.2307  # H3: 2225.3635 -> 2398.2307
$ws.Cells.Item(3, 9).Value = 2017  # I3: 1721 -> 2017
$ws.Cells.Item(3, 11).Value = 2017  # K3: 1721 -> 2017
$ws.Cells.Item(3, 13).Value = -1903  # M3: -1607 -> -1903

$ws.Cells.Item(60, 8).Value = 0  # H60: 92949.5 -> 0
$ws.Cells.Item(60, 10).Value = 0  # J60: 92949.5 -> 0
$ws.Cells.Item(60, 12).Value = 0  # L60: 92949.5 -> 0
$ws.Cells.Item(60, 14).ClearContents()  # remove N60

$ws.Cells.Item(105, 8).Value = 2001.25  # H105: 2005 -> 2001.25
$ws.Cells.Item(105, 9).Value = 2001.25  # I105: 2005 -> 2001.25
$ws.Cells.Item(105, 11).Value = 2001.25  # K105: 2005 -> 2001.25
$ws.Cells.Item(105, 13).Value = -254.25  # M105: -258 -> -254.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 67.14286  # H7: 67.875 -> 67.14286
$ws.Cells.Item(7, 10).Value = 0  # J7: 73 -> 0
$ws.Cells.Item(7, 12).Value = 0  # L7: 73 -> 0
$ws.Cells.Item(7, 14).ClearContents()  # remove N7

$ws.Cells.Item(11, 8).Value = 1264  # H11: 847.5 -> 1264
$ws.Cells.Item(11, 9).Value = 1350  # I11: 1522.5 -> 1350
$ws.Cells.Item(11, 10).Value = 1006  # J11: 510 -> 1006
$ws.Cells.Item(11, 11).Value = 1350  # K11: 1522.5 -> 1350
$ws.Cells.Item(11, 12).Value = 1006  # L11: 510 -> 1006
$ws.Cells.Item(11, 13).Value = -1210  # M11: -1382.5 -> -1210
$ws.Cells.Item(11, 14).Value = -1286  # N11: -790 -> -1286

$ws.Cells.Item(15, 8).Value = 9999  # H15: 6724.6665 -> 9999
$ws.Cells.Item(15, 9).Value = 0  # I15: 175 -> 0
$ws.Cells.Item(15, 10).Value = 9999  # J15: 9999.5 -> 9999
$ws.Cells.Item(15, 11).Value = 0  # K15: 175 -> 0
$ws.Cells.Item(15, 12).Value = 9999  # L15: 9999.5 -> 9999
$ws.Cells.Item(15, 14).Value = -10339  # N15: -10339.5 -> -10339
$ws.Cells.Item(15, 13).ClearContents()  # remove M15

$ws.Cells.Item(31, 8).Value = 8304  # H31: 6376 -> 8304
$ws.Cells.Item(31, 9).Value = 0  # I31: 1000 -> 0
$ws.Cells.Item(31, 10).Value = 8304  # J31: 7988.8 -> 8304
$ws.Cells.Item(31, 11).Value = 0  # K31: 1000 -> 0
$ws.Cells.Item(31, 12).Value = 8304  # L31: 7988.8 -> 8304
$ws.Cells.Item(31, 14).Value = -8894  # N31: -8578.799999999999 -> -8894
$ws.Cells.Item(31, 13).ClearContents()  # remove M31

$ws.Cells.Item(34, 8).Value = 8304  # H34: 6376 -> 8304
$ws.Cells.Item(34, 9).Value = 0  # I34: 1000 -> 0
$ws.Cells.Item(34, 10).Value = 8304  # J34: 7988.8 -> 8304
$ws.Cells.Item(34, 11).Value = 0  # K34: 1000 -> 0
$ws.Cells.Item(34, 12).Value = 8304  # L34: 7988.8 -> 8304
$ws.Cells.Item(34, 14).Value = -8708  # N34: -8392.799999999999 -> -8708
$ws.Cells.Item(34, 13).ClearContents()  # remove M34

$ws.Cells.Item(107, 8).Value = 299.7647  # H107: 347.75 -> 299.7647
$ws.Cells.Item(107, 9).Value = 379.2  # I107: 412.33334 -> 379.2
$ws.Cells.Item(107, 10).Value = 186.28572  # J107: 264.7143 -> 186.28572
$ws.Cells.Item(107, 11).Value = 379.2  # K107: 412.33334 -> 379.2
$ws.Cells.Item(107, 12).Value = 186.28572  # L107: 264.7143 -> 186.28572
$ws.Cells.Item(107, 13).Value = 1540.8  # M107: 1507.66666 -> 1540.8
$ws.Cells.Item(107, 14).Value = -4026.28572  # N107: -4104.7143 -> -4026.28572

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(80, 8).Value = 4905.4165  # H80: 4979.696 -> 4905.4165
$ws.Cells.Item(80, 9).Value = 4899.8335  # I80: 5000 -> 4899.8335
$ws.Cells.Item(80, 11).Value = 14699.5005  # K80: 15000 -> 14699.5005
$ws.Cells.Item(80, 13).Value = -13763.5005  # M80: -14064 -> -13763.5005

$ws.Cells.Item(83, 8).Value = 4905.4165  # H83: 4979.696 -> 4905.4165
$ws.Cells.Item(83, 9).Value = 4899.8335  # I83: 5000 -> 4899.8335
$ws.Cells.Item(83, 11).Value = 44098.5015  # K83: 45000 -> 44098.5015
$ws.Cells.Item(83, 13).Value = -39418.5015  # M83: -40320 -> -39418.5015

$ws.Cells.Item(105, 8).Value = 30000  # H105: 0 -> 30000
$ws.Cells.Item(105, 9).Value = 30000  # I105: 0 -> 30000
$ws.Cells.Item(105, 11).Value = 90000  # K105: 0 -> 90000
$ws.Cells.Item(105, 13).Value = -87379  # add M105 = -87379

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5555  # H70: 0 -> 5555
$ws.Cells.Item(70, 9).Value = 5555  # I70: 0 -> 5555
$ws.Cells.Item(70, 11).Value = 5555  # K70: 0 -> 5555
$ws.Cells.Item(70, 13).Value = -5285  # add M70 = -5285

$ws.Cells.Item(73, 8).Value = 5555  # H73: 0 -> 5555
$ws.Cells.Item(73, 9).Value = 5555  # I73: 0 -> 5555
$ws.Cells.Item(73, 11).Value = 5555  # K73: 0 -> 5555
$ws.Cells.Item(73, 13).Value = -4619  # add M73 = -4619

$ws.Cells.Item(80, 8).Value = 2347.8333  # H80: 2467.4 -> 2347.8333
$ws.Cells.Item(80, 9).Value = 1749.6666  # I80: 1749.5 -> 1749.6666
$ws.Cells.Item(80, 11).Value = 1749.6666  # K80: 1749.5 -> 1749.6666
$ws.Cells.Item(80, 13).Value = -751.6666  # M80: -751.5 -> -751.6666

$ws.Cells.Item(83, 8).Value = 2347.8333  # H83: 2467.4 -> 2347.8333
$ws.Cells.Item(83, 9).Value = 1749.6666  # I83: 1749.5 -> 1749.6666
$ws.Cells.Item(83, 11).Value = 8748.333000000001  # K83: 8747.5 -> 8748.333000000001
$ws.Cells.Item(83, 13).Value = -3756.333000000001  # M83: -3755.5 -> -3756.333000000001

$ws.Cells.Item(97, 8).Value = 505.78262  # H97: 507.30435 -> 505.78262
$ws.Cells.Item(97, 9).Value = 505.78262  # I97: 527.1818 -> 505.78262
$ws.Cells.Item(97, 10).Value = 0  # J97: 70 -> 0
$ws.Cells.Item(97, 11).Value = 505.78262  # K97: 527.1818 -> 505.78262
$ws.Cells.Item(97, 12).Value = 0  # L97: 70 -> 0
$ws.Cells.Item(97, 13).Value = -9.782620000000009  # M97: -31.18179999999995 -> -9.782620000000009
$ws.Cells.Item(97, 14).ClearContents()  # remove N97

$ws.Cells.Item(132, 8).Value = 13542.667  # H132: 14169.875 -> 13542.667
$ws.Cells.Item(132, 9).Value = 14206.6  # I132: 15049.923 -> 14206.6
$ws.Cells.Item(132, 10).Value = 10223  # J132: 10356.333 -> 10223
$ws.Cells.Item(132, 11).Value = 42619.8  # K132: 45149.769 -> 42619.8
$ws.Cells.Item(132, 12).Value = 30669  # L132: 31068.999 -> 30669
$ws.Cells.Item(132, 13).Value = -40089.8  # M132: -42619.769 -> -40089.8
$ws.Cells.Item(132, 14).Value = -35729  # N132: -36128.999 -> -35729

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1000.25  # H16: 999.6667 -> 1000.25

$ws.Cells.Item(40, 8).Value = 5908.9165  # H40: 6319.8096 -> 5908.9165
$ws.Cells.Item(40, 9).Value = 5264.5625  # I40: 5779.615 -> 5264.5625
$ws.Cells.Item(40, 11).Value = 5264.5625  # K40: 5779.615 -> 5264.5625
$ws.Cells.Item(40, 13).Value = -5128.5625  # M40: -5643.615 -> -5128.5625

$ws.Cells.Item(55, 8).Value = 1126.9  # H55: 1316.375 -> 1126.9
$ws.Cells.Item(55, 9).Value = 1833.8  # I55: 1588.5 -> 1833.8
$ws.Cells.Item(55, 10).Value = 420  # J55: 500 -> 420
$ws.Cells.Item(55, 11).Value = 1833.8  # K55: 1588.5 -> 1833.8
$ws.Cells.Item(55, 12).Value = 420  # L55: 500 -> 420
$ws.Cells.Item(55, 13).Value = -1660.8  # M55: -1415.5 -> -1660.8
$ws.Cells.Item(55, 14).Value = -766  # N55: -846 -> -766

$ws.Cells.Item(68, 8).Value = 7611.615  # H68: 7912.5835 -> 7611.615
$ws.Cells.Item(68, 10).Value = 9157.286  # J68: 10016.833 -> 9157.286
$ws.Cells.Item(68, 12).Value = 9157.286  # L68: 10016.833 -> 9157.286
$ws.Cells.Item(68, 14).Value = -10655.286  # N68: -11514.833 -> -10655.286

$ws.Cells.Item(71, 8).Value = 7611.615  # H71: 7912.5835 -> 7611.615
$ws.Cells.Item(71, 10).Value = 9157.286  # J71: 10016.833 -> 9157.286
$ws.Cells.Item(71, 12).Value = 45786.43  # L71: 50084.165 -> 45786.43
$ws.Cells.Item(71, 14).Value = -53274.43  # N71: -57572.165 -> -53274.43

$ws.Cells.Item(93, 8).Value = 983.1579  # H93: 1024.5 -> 983.1579
$ws.Cells.Item(93, 9).Value = 983.1579  # I93: 1024.5 -> 983.1579
$ws.Cells.Item(93, 11).Value = 983.1579  # K93: 1024.5 -> 983.1579
$ws.Cells.Item(93, 13).Value = 264.8421  # M93: 223.5 -> 264.8421

$ws.Cells.Item(122, 8).Value = 4129.3335  # H122: 3672 -> 4129.3335
$ws.Cells.Item(122, 9).Value = 3595.182  # I122: 3672 -> 3595.182
$ws.Cells.Item(122, 10).Value = 10005  # J122: 0 -> 10005
$ws.Cells.Item(122, 11).Value = 10785.546  # K122: 11016 -> 10785.546
$ws.Cells.Item(122, 12).Value = 30015  # L122: 0 -> 30015
$ws.Cells.Item(122, 13).Value = -8335.545999999998  # M122: -8566 -> -8335.545999999998
$ws.Cells.Item(122, 14).Value = -34915  # add N122 = -34915

$ws.Cells.Item(136, 8).Value = 3556.1667  # H136: 3598.6667 -> 3556.1667
$ws.Cells.Item(136, 9).Value = 3071.5  # I136: 3122.5 -> 3071.5
$ws.Cells.Item(136, 11).Value = 9214.5  # K136: 9367.5 -> 9214.5
$ws.Cells.Item(136, 13).Value = -6664.5  # M136: -6817.5 -> -6664.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 1343  # H122: 1385.875 -> 1343
$ws.Cells.Item(122, 9).Value = 1343  # I122: 1385.875 -> 1343
$ws.Cells.Item(122, 11).Value = 4029  # K122: 4157.625 -> 4029
$ws.Cells.Item(122, 13).Value = -1579  # M122: -1707.625 -> -1579

$ws.Cells.Item(126, 8).Value = 3614.4443  # H126: 3074.5217 -> 3614.4443
$ws.Cells.Item(126, 9).Value = 1462.9166  # I126: 1396.3572 -> 1462.9166
$ws.Cells.Item(126, 10).Value = 7917.5  # J126: 5685 -> 7917.5
$ws.Cells.Item(126, 11).Value = 4388.7498  # K126: 4189.071599999999 -> 4388.7498
$ws.Cells.Item(126, 12).Value = 23752.5  # L126: 17055 -> 23752.5
$ws.Cells.Item(126, 13).Value = -1918.7498  # M126: -1719.071599999999 -> -1918.7498
$ws.Cells.Item(126, 14).Value = -28692.5  # N126: -21995 -> -28692.5

$ws.Cells.Item(132, 8).Value = 3201.5454  # H132: 3268.5 -> 3201.5454
$ws.Cells.Item(132, 9).Value = 3201.5454  # I132: 3202 -> 3201.5454
$ws.Cells.Item(132, 10).Value = 0  # J132: 4000 -> 0
$ws.Cells.Item(132, 11).Value = 9604.636200000001  # K132: 9606 -> 9604.636200000001
$ws.Cells.Item(132, 12).Value = 0  # L132: 12000 -> 0
$ws.Cells.Item(132, 13).Value = -7074.636200000001  # M132: -7076 -> -7074.636200000001
$ws.Cells.Item(132, 14).ClearContents()  # remove N132
